$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.01829118917108
$ws.Cells.Item(2, 4).Value = 1.023985953754433
$ws.Cells.Item(2, 5).Value = 1.019583685794839
$ws.Cells.Item(2, 9).Value = 1.027341434458993
$ws.Cells.Item(2, 10).Value = 1.023500457630093
$ws.Cells.Item(2, 11).Value = 1.026815782312966
$ws.Cells.Item(2, 12).Value = 1.022426501826416
$ws.Cells.Item(2, 14).Value = 1.024953945677072

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.01921308746259
$ws.Cells.Item(3, 4).Value = 1.024666311906729
$ws.Cells.Item(3, 5).Value = 1.020363800528458
$ws.Cells.Item(3, 9).Value = 1.027486096304845
$ws.Cells.Item(3, 10).Value = 1.024058758255928
$ws.Cells.Item(3, 11).Value = 1.027303480658399
$ws.Cells.Item(3, 12).Value = 1.023012744025828
$ws.Cells.Item(3, 14).Value = 1.025513039153833

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.019809471143432
$ws.Cells.Item(4, 4).Value = 1.025105815754036
$ws.Cells.Item(4, 5).Value = 1.020868883288441
$ws.Cells.Item(4, 9).Value = 1.027577567650423
$ws.Cells.Item(4, 10).Value = 1.024419280685079
$ws.Cells.Item(4, 11).Value = 1.027617646555353
$ws.Cells.Item(4, 12).Value = 1.023391725065711
$ws.Cells.Item(4, 14).Value = 1.025874073566186

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.020060154882545
$ws.Cells.Item(5, 4).Value = 1.025290406023527
$ws.Cells.Item(5, 5).Value = 1.021081289651895
$ws.Cells.Item(5, 9).Value = 1.027615510338181
$ws.Cells.Item(5, 10).Value = 1.024570667198113
$ws.Cells.Item(5, 11).Value = 1.027749383663829
$ws.Cells.Item(5, 12).Value = 1.023550962170415
$ws.Cells.Item(5, 14).Value = 1.026025675065433

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.020102243648351
$ws.Cells.Item(6, 4).Value = 1.025321389107824
$ws.Cells.Item(6, 5).Value = 1.021116957654372
$ws.Cells.Item(6, 9).Value = 1.02762185104137
$ws.Cells.Item(6, 10).Value = 1.024596075242659
$ws.Cells.Item(6, 11).Value = 1.02777148303271
$ws.Cells.Item(6, 12).Value = 1.023577693675258
$ws.Cells.Item(6, 14).Value = 1.026051119192316

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.019812820934355
$ws.Cells.Item(7, 4).Value = 1.025108282956336
$ws.Cells.Item(7, 5).Value = 1.02087172120105
$ws.Cells.Item(7, 9).Value = 1.027578076654849
$ws.Cells.Item(7, 10).Value = 1.024421304216371
$ws.Cells.Item(7, 11).Value = 1.027619408163696
$ws.Cells.Item(7, 12).Value = 1.023393853141069
$ws.Cells.Item(7, 14).Value = 1.025876099971125

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.018602778697911
$ws.Cells.Item(8, 4).Value = 1.024216034883061
$ws.Cells.Item(8, 5).Value = 1.019847267031813
$ws.Cells.Item(8, 9).Value = 1.027390765068835
$ws.Cells.Item(8, 10).Value = 1.023689289578679
$ws.Cells.Item(8, 11).Value = 1.026980892936636
$ws.Cells.Item(8, 12).Value = 1.022624698360698
$ws.Cells.Item(8, 14).Value = 1.025143045788683

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.016469452398451
$ws.Cells.Item(9, 4).Value = 1.022638231884035
$ws.Cells.Item(9, 5).Value = 1.018044378008663
$ws.Cells.Item(9, 9).Value = 1.027044387154979
$ws.Cells.Item(9, 10).Value = 1.022393800020861
$ws.Cells.Item(9, 11).Value = 1.025845022556851
$ws.Cells.Item(9, 12).Value = 1.021266663828438
$ws.Cells.Item(9, 14).Value = 1.023845716487099

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.015046564942677
$ws.Cells.Item(10, 4).Value = 1.021582724231355
$ws.Cells.Item(10, 5).Value = 1.016844098516479
$ws.Cells.Item(10, 9).Value = 1.026802550850059
$ws.Cells.Item(10, 10).Value = 1.021526445035799
$ws.Cells.Item(10, 11).Value = 1.025080637452521
$ws.Cells.Item(10, 12).Value = 1.020359564662427
$ws.Cells.Item(10, 14).Value = 1.022977129758471

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.014430290770159
$ws.Cells.Item(11, 4).Value = 1.021124833331282
$ws.Cells.Item(11, 5).Value = 1.016324771183726
$ws.Cells.Item(11, 9).Value = 1.026695253481622
$ws.Cells.Item(11, 10).Value = 1.02115000656751
$ws.Cells.Item(11, 11).Value = 1.024747971992795
$ws.Cells.Item(11, 12).Value = 1.019966380145436
$ws.Cells.Item(11, 14).Value = 1.022600156704379

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.014201356663806
$ws.Cells.Item(12, 4).Value = 1.020954626140263
$ws.Cells.Item(12, 5).Value = 1.016131931506012
$ws.Cells.Item(12, 9).Value = 1.026655011494061
$ws.Cells.Item(12, 10).Value = 1.021010051078788
$ws.Cells.Item(12, 11).Value = 1.024624153706837
$ws.Cells.Item(12, 12).Value = 1.019820274107633
$ws.Cells.Item(12, 14).Value = 1.022460002462809

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.014250464809797
$ws.Cells.Item(13, 4).Value = 1.020991141844759
$ws.Cells.Item(13, 5).Value = 1.016173293462756
$ws.Cells.Item(13, 9).Value = 1.026663661036206
$ws.Cells.Item(13, 10).Value = 1.021040077850614
$ws.Cells.Item(13, 11).Value = 1.024650724513404
$ws.Cells.Item(13, 12).Value = 1.019851617039155
$ws.Cells.Item(13, 14).Value = 1.022490071876096

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.014411367458804
$ws.Cells.Item(14, 4).Value = 1.021110766516329
$ws.Cells.Item(14, 5).Value = 1.01630882973148
$ws.Cells.Item(14, 9).Value = 1.026691934956114
$ws.Cells.Item(14, 10).Value = 1.021138440439163
$ws.Cells.Item(14, 11).Value = 1.024737742259148
$ws.Cells.Item(14, 12).Value = 1.01995430418905
$ws.Cells.Item(14, 14).Value = 1.022588574150803

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.014510501940987
$ws.Cells.Item(15, 4).Value = 1.021184454555114
$ws.Cells.Item(15, 5).Value = 1.016392346302831
$ws.Cells.Item(15, 9).Value = 1.026709304196541
$ws.Cells.Item(15, 10).Value = 1.021199027746615
$ws.Cells.Item(15, 11).Value = 1.024791323468656
$ws.Cells.Item(15, 12).Value = 1.020017565236184
$ws.Cells.Item(15, 14).Value = 1.022649247499182

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.01508746184072
$ws.Cells.Item(16, 4).Value = 1.02161309517626
$ws.Cells.Item(16, 5).Value = 1.016878573132054
$ws.Cells.Item(16, 9).Value = 1.026809617515564
$ws.Cells.Item(16, 10).Value = 1.021551409808405
$ws.Cells.Item(16, 11).Value = 1.02510268002937
$ws.Cells.Item(16, 12).Value = 1.020385650583772
$ws.Cells.Item(16, 14).Value = 1.023002129983917

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.015449332724244
$ws.Cells.Item(17, 4).Value = 1.021881743806143
$ws.Cells.Item(17, 5).Value = 1.017183678804373
$ws.Cells.Item(17, 9).Value = 1.026871851026292
$ws.Cells.Item(17, 10).Value = 1.021772218028952
$ws.Cells.Item(17, 11).Value = 1.025297536345466
$ws.Cells.Item(17, 12).Value = 1.020616433174083
$ws.Cells.Item(17, 14).Value = 1.023223251777465

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.015660390928257
$ws.Cells.Item(18, 4).Value = 1.022038359982619
$ws.Cells.Item(18, 5).Value = 1.017361680445129
$ws.Cells.Item(18, 9).Value = 1.026907901735118
$ws.Cells.Item(18, 10).Value = 1.021900927913351
$ws.Cells.Item(18, 11).Value = 1.025411030406618
$ws.Cells.Item(18, 12).Value = 1.02075100565172
$ws.Cells.Item(18, 14).Value = 1.023352144444663

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.015732353760437
$ws.Cells.Item(19, 4).Value = 1.022091748094493
$ws.Cells.Item(19, 5).Value = 1.017422380923675
$ws.Cells.Item(19, 9).Value = 1.026920151829983
$ws.Cells.Item(19, 10).Value = 1.021944800412823
$ws.Cells.Item(19, 11).Value = 1.02544970137978
$ws.Cells.Item(19, 12).Value = 1.020796884726364
$ws.Cells.Item(19, 14).Value = 1.023396079248116

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.01541050894483
$ws.Cells.Item(20, 4).Value = 1.021852928824224
$ws.Cells.Item(20, 5).Value = 1.017150939844709
$ws.Cells.Item(20, 9).Value = 1.026865199714017
$ws.Cells.Item(20, 10).Value = 1.021748536059155
$ws.Cells.Item(20, 11).Value = 1.025276646889304
$ws.Cells.Item(20, 12).Value = 1.020591676423659
$ws.Cells.Item(20, 14).Value = 1.023199536176554

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.014363986252301
$ws.Cells.Item(21, 4).Value = 1.021075543502177
$ws.Cells.Item(21, 5).Value = 1.016268915967339
$ws.Cells.Item(21, 9).Value = 1.026683619667659
$ws.Cells.Item(21, 10).Value = 1.021109478675319
$ws.Cells.Item(21, 11).Value = 1.024712124636128
$ws.Cells.Item(21, 12).Value = 1.019924067025775
$ws.Cells.Item(21, 14).Value = 1.022559571257933

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.013705866452491
$ws.Cells.Item(22, 4).Value = 1.020586041492744
$ws.Cells.Item(22, 5).Value = 1.015714709670309
$ws.Cells.Item(22, 9).Value = 1.026567214573774
$ws.Cells.Item(22, 10).Value = 1.020706929253519
$ws.Cells.Item(22, 11).Value = 1.024355732409183
$ws.Cells.Item(22, 12).Value = 1.019503968986966
$ws.Cells.Item(22, 14).Value = 1.02215645016978

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.014054760223465
$ws.Cells.Item(23, 4).Value = 1.02084560439103
$ws.Cells.Item(23, 5).Value = 1.016008470638415
$ws.Cells.Item(23, 9).Value = 1.026629135054783
$ws.Cells.Item(23, 10).Value = 1.020920398975989
$ws.Cells.Item(23, 11).Value = 1.024544800185112
$ws.Cells.Item(23, 12).Value = 1.019726703330257
$ws.Cells.Item(23, 14).Value = 1.022370223043741

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.015428051785808
$ws.Cells.Item(24, 4).Value = 1.021865949327775
$ws.Cells.Item(24, 5).Value = 1.017165733051983
$ws.Cells.Item(24, 9).Value = 1.026868205925281
$ws.Cells.Item(24, 10).Value = 1.021759237180948
$ws.Cells.Item(24, 11).Value = 1.02528608643701
$ws.Cells.Item(24, 12).Value = 1.020602863054416
$ws.Cells.Item(24, 14).Value = 1.023210252495167

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.017021090411649
$ws.Cells.Item(25, 4).Value = 1.023046779788629
$ws.Cells.Item(25, 5).Value = 1.018510183636503
$ws.Cells.Item(25, 9).Value = 1.027135862007552
$ws.Cells.Item(25, 10).Value = 1.022729370841061
$ws.Cells.Item(25, 11).Value = 1.026139935700584
$ws.Cells.Item(25, 12).Value = 1.021618060051877
$ws.Cells.Item(25, 14).Value = 1.024181763856354

